$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "Summary" section header row (old row 39); this shifts
# the six "Total ..." rows up by one (old 40-45 -> new 39-44).
$ws.Rows("39").Delete()

# Prefix each category's detail-row labels with the category name, e.g.
# "     New nominations" -> "     Civilian, New nominations".
$ws.Range("A7").Value = "     Civilian, New nominations"
$ws.Range("A8").Value = "     Civilian, Carryover nominations"
$ws.Range("A9").Value = "     Civilian, Confirmed "
$ws.Range("A10").Value = "     Civilian, Withdrawn "
$ws.Range("A11").Value = "     Civilian, Returned to White House "

$ws.Range("A13").Value = "     Other Civilian, New nominations"
$ws.Range("A14").Value = "     Other Civilian, Carryover nominations"
$ws.Range("A15").Value = "     Other Civilian, Confirmed "
$ws.Range("A16").Value = "     Other Civilian, Withdrawn "
$ws.Range("A17").Value = "     Other Civilian, Returned to White House "

$ws.Range("A19").Value = "     Air Force, New nominations"
$ws.Range("A20").Value = "     Air Force, Carryover nominations"
$ws.Range("A21").Value = "     Air Force, Confirmed "
$ws.Range("A22").Value = "     Air Force, Returned to White House "

$ws.Range("A24").Value = "     Army, New nominations"
$ws.Range("A25").Value = "     Army, Carryover nominations"
$ws.Range("A26").Value = "     Army, Confirmed "
$ws.Range("A27").Value = "     Army, Returned to White House "

$ws.Range("A29").Value = "     Navy, New nominations"
$ws.Range("A30").Value = "     Navy, Carryover nominations"
$ws.Range("A31").Value = "     Navy, Confirmed "
$ws.Range("A32").Value = "     Navy, Withdrawn "
$ws.Range("A33").Value = "     Navy, Returned to White House "

$ws.Range("A35").Value = "     Marine Corps, New nominations"
$ws.Range("A36").Value = "     Marine Corps, Carryover nominations"
$ws.Range("A37").Value = "     Marine Corps, Confirmed "
$ws.Range("A38").Value = "     Marine Corps, Returned to White House "

# The old Summary block listed "Total nominations carried over from the
# First Session" (2207) then "Total nominations received this Session"
# (19680). The new layout swaps the order and renames them to "Total new
# nominations" (19680) and "Total carryover nominations" (2207).
$ws.Range("A39").Value = "Total new nominations"
$ws.Range("B39").Value = 19680
$ws.Range("A40").Value = "Total carryover nominations"
$ws.Range("B40").Value = 2207
